$d = $word.ActiveDocument
Write-Host "Paragraphs count:" $d.Paragraphs.Count
Write-Host "Content length:" $d.Content.Text.Length
